$d = $word.ActiveDocument
$d.Content.Find.Execute("99÷9=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "61÷9=6, 7", 2) | Out-Null
$d.Content.Find.Execute("76÷6=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=4, 7", 2) | Out-Null
$d.Content.Find.Execute("42÷2=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "47÷7=6, 5", 2) | Out-Null
$d.Content.Find.Execute("95÷5=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2) | Out-Null
$d.Content.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "85÷3=28, 1", 2) | Out-Null
$d.Content.Find.Execute("49÷8=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=46, 1", 2) | Out-Null
$d.Content.Find.Execute("70÷7=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷9=10, 2", 2) | Out-Null
$d.Content.Find.Execute("80÷5=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=2, 4", 2) | Out-Null
$d.Content.Find.Execute("60÷7=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=1, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷4=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "51÷2=25, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "47÷3=15, 2", 2) | Out-Null
$d.Content.Find.Execute("61÷3=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "60÷9=6, 6", 2) | Out-Null
$d.Content.Find.Execute("83÷9=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "57÷4=14, 1", 2) | Out-Null
$d.Content.Find.Execute("60÷8=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "43÷5=8, 3", 2) | Out-Null
$d.Content.Find.Execute("98÷3=32, 2", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=6, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷2=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2) | Out-Null
$d.Content.Find.Execute("18÷4=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=11, 0", 2) | Out-Null
$d.Content.Find.Execute("92÷8=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "36÷7=5, 1", 2) | Out-Null
$d.Content.Find.Execute("62÷9=6, 8", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=9, 4", 2) | Out-Null
$d.Content.Find.Execute("10÷8=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2) | Out-Null
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=25, 0", 2) | Out-Null
$d.Content.Find.Execute("70÷6=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "63÷8=7, 7", 2) | Out-Null
$d.Content.Find.Execute("15÷5=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=4, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷4=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=18, 1", 2) | Out-Null
$d.Content.Find.Execute("97÷5=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "83÷5=16, 3", 2) | Out-Null
